$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the new "Числ. врачей" value under the Здравоохранение column (E20 was blank)
$ws.Range("E20").Value = "Числ. врачей - docsnum (чел.) (8018106)"

# 2. Add a new category block (rows 24-26), copying the look (styles/borders/fill)
#    of an existing "категория/признаки" block (rows 10-12) that already has the
#    B:s=1/2/4 + C:s=1/2/2 styling we need, then overwrite with the right labels.
$ws.Range("B10:C12").Copy() | Out-Null
$ws.Range("B24:C26").PasteSpecial(-4122) | Out-Null

$ws.Range("B24").Value = "категория"
$ws.Range("B25").Value = "признаки"

# 3. Move the active-cell selection from E28 to D28 (no data there, just UI state)
$ws.Range("D28").Select() | Out-Null
